# LLD and HLD Added
# Removes stray/empty placeholder shapes left over from slide-layout
# rearrangement: the empty "Title 1" (ctrTitle) placeholder on the title
# slide, and the empty "Content Placeholder 2" placeholders on the
# content slides.

$p = $ppt.ActivePresentation

# Slide 1: remove the empty ctrTitle placeholder ("Title 1"), keep Subtitle 2
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("Title 1").Delete()

# Slides 2, 3, 5, 6, 7, 8: remove the empty "Content Placeholder 2" shape
$targetSlides = @(2, 3, 5, 6, 7, 8)
foreach ($idx in $targetSlides) {
    $slide = $p.Slides.Item($idx)
    $slide.Shapes.Item("Content Placeholder 2").Delete()
}
